$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
WITH diagnosis_summary AS ( 
    SELECT 
        dgn."participant.id" AS participant_id, 
        GROUP_CONCAT(dgn.age_at_diagnosis) AS age_at_diagnosis,  -- Removed DISTINCT from GROUP_CONCAT
        GROUP_CONCAT(dgn.diagnosis) AS diagnosis,                 -- Removed DISTINCT from GROUP_CONCAT
        GROUP_CONCAT(dgn.anatomic_site) AS anatomic_site,         -- Removed DISTINCT from GROUP_CONCAT
        GROUP_CONCAT(dgn.diagnosis_basis) AS diagnosis_basis      -- Removed DISTINCT from GROUP_CONCAT
    FROM df_diagnosis dgn
    WHERE dgn."participant.id" IS NOT NULL
    GROUP BY dgn."participant.id"
),
survival_summary AS (
    SELECT 
        "participant.id" AS participant_id,
        GROUP_CONCAT(last_known_survival_status) AS last_known_survival_status  -- Removed DISTINCT from GROUP_CONCAT
    FROM df_survival
    WHERE last_known_survival_status IS NOT NULL
    GROUP BY "participant.id"
),
treatment_summary AS (
    SELECT 
        trt."participant.id" AS participant_id,
        GROUP_CONCAT(trt.treatment_type) AS treatment_type  -- Removed DISTINCT from GROUP_CONCAT
    FROM df_treatment trt
    WHERE trt."participant.id" IS NOT NULL
    GROUP BY trt."participant.id"
)

SELECT
    prt.participant_id AS "Participant ID",
    std.dbgap_accession AS "Study ID",
    COALESCE(prt.sex_at_birth, '') AS "Sex",
    COALESCE(prt.race, '') AS "Race",
    COALESCE(dgn.diagnosis, '') AS "Diagnosis",
    COALESCE(dgn.anatomic_site, '') AS "Diagnosis Anatomic Site",
    COALESCE(
        CASE 
            WHEN dgn.age_at_diagnosis = '-999' THEN 'Not Reported'
            ELSE dgn.age_at_diagnosis 
        END, ''
    ) AS "Age at Diagnosis (days)",
    COALESCE(trt.treatment_type, '') AS "Treatment Type",
    COALESCE(srv.last_known_survival_status, '') AS "Last Known Survival Status"
FROM df_participant prt
JOIN df_study std ON prt."study.id" = std.id
LEFT JOIN diagnosis_summary dgn ON prt.id = dgn.participant_id
LEFT JOIN survival_summary srv ON prt.id = srv.participant_id
LEFT JOIN treatment_summary trt ON prt.id = trt.participant_id
WHERE std.dbgap_accession = 'phs002504'
  AND EXISTS (
      SELECT 1 
      FROM df_sample smp
      JOIN df_sequencing_file f ON f."sample.id" = smp.id
      WHERE smp."participant.id" = prt.id
        AND f.file_type = 'bam'
        AND f.library_selection = 'PCR'
        AND f.library_strategy = 'WXS'
        AND f.data_category = 'Sequencing'
  )
ORDER BY prt.participant_id ASC
LIMIT 100;
'@

$ws.Range("B2").Value = $newQuery

$ws.Range("B2").Select() | Out-Null
